$wb = $excel.ActiveWorkbook

# Sheets (in tab order): 1=userNamePage, 2=passwordPage, 3=homePageValidation, 4=productSearchPage
$wsPassword = $wb.Worksheets.Item(2)
$wsProductSearch = $wb.Worksheets.Item(4)

# --- productSearchPage: move the two data cells from column D to column A ---
$d1Value = $wsProductSearch.Range("D1").Value2
$d2Value = $wsProductSearch.Range("D2").Value2

$wsProductSearch.Range("A1").Value = $d1Value
$wsProductSearch.Range("A2").Value = $d2Value

[void]$wsProductSearch.Range("D1").ClearContents()
[void]$wsProductSearch.Range("D2").ClearContents()
[void]$wsProductSearch.Range("A2").ClearFormats()

# --- Selections on each sheet ---
[void]$wsPassword.Range("A3").Select()

# Activate productSearchPage last so it becomes the active/selected tab
$wsProductSearch.Select()
[void]$wsProductSearch.Range("A2").Select()
